# Updated symbol list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'322.47"
$ws.Cells.Item(2, 5).Value = "'-2.75%"
$ws.Cells.Item(3, 4).Value = "'42.96"
$ws.Cells.Item(3, 5).Value = "'-5.52%"
$ws.Cells.Item(4, 4).Value = "'5.194"
$ws.Cells.Item(4, 5).Value = "'-7.81%"
$ws.Cells.Item(5, 4).Value = "'0.08200"
$ws.Cells.Item(6, 4).Value = "'4.324"
$ws.Cells.Item(6, 5).Value = "'-2.68%"
$ws.Cells.Item(7, 4).Value = "'1.828"
$ws.Cells.Item(7, 5).Value = "'-11.40%"
$ws.Cells.Item(8, 4).Value = "'0.9341"
$ws.Cells.Item(8, 5).Value = "'-3.49%"
$ws.Cells.Item(9, 4).Value = "'0.1113"
$ws.Cells.Item(9, 5).Value = "'-4.97%"
$ws.Cells.Item(10, 4).Value = "'0.1861"
$ws.Cells.Item(10, 5).Value = "'-3.03%"
$ws.Cells.Item(11, 4).Value = "'0.09508"
$ws.Cells.Item(11, 5).Value = "'-3.59%"
$ws.Cells.Item(12, 4).Value = "'0.04608"
$ws.Cells.Item(12, 5).Value = "'-0.48%"
$ws.Cells.Item(13, 4).Value = "'7.420"
$ws.Cells.Item(13, 5).Value = "'-28.40%"
$ws.Cells.Item(14, 4).Value = "'0.1056"
$ws.Cells.Item(14, 5).Value = "'-0.29%"
$ws.Cells.Item(15, 4).Value = "'0.001300"
$ws.Cells.Item(15, 5).Value = "'1.21%"
$ws.Cells.Item(16, 4).Value = "'0.005904"
$ws.Cells.Item(16, 5).Value = "'-3.31%"
$ws.Cells.Item(17, 4).Value = "'3.363"
$ws.Cells.Item(17, 5).Value = "'-0.41%"
$ws.Cells.Item(18, 4).Value = "'2.546"
$ws.Cells.Item(18, 5).Value = "'0.08%"
$ws.Cells.Item(19, 4).Value = "'0.3369"
$ws.Cells.Item(19, 5).Value = "'0.01%"
$ws.Cells.Item(20, 4).Value = "'0.1387"
$ws.Cells.Item(20, 5).Value = "'-0.48%"
$ws.Cells.Item(21, 4).Value = "'0.2619"
$ws.Cells.Item(21, 5).Value = "'-1.37%"
$ws.Cells.Item(22, 4).Value = "'0.04161"
$ws.Cells.Item(22, 5).Value = "'-0.90%"
$ws.Cells.Item(23, 4).Value = "'0.001242"
$ws.Cells.Item(23, 5).Value = "'-5.40%"
$ws.Cells.Item(24, 4).Value = "'0.004336"
$ws.Cells.Item(24, 5).Value = "'-4.94%"
$ws.Cells.Item(25, 4).Value = "'0.0001198"
$ws.Cells.Item(25, 5).Value = "'-8.17%"
$ws.Cells.Item(26, 4).Value = "'0.0002974"
$ws.Cells.Item(26, 5).Value = "'-20.70%"
$ws.Cells.Item(38, 4).Value = "'0.02771"
$ws.Cells.Item(38, 5).Value = "'1.05%"
$ws.Cells.Item(39, 4).Value = "'0.05584"
$ws.Cells.Item(39, 5).Value = "'-3.14%"
$ws.Cells.Item(40, 4).Value = "'0.008056"
$ws.Cells.Item(40, 5).Value = "'4.91%"
$ws.Cells.Item(41, 4).Value = "'0.1397"
$ws.Cells.Item(41, 5).Value = "'-2.71%"
$ws.Cells.Item(42, 4).Value = "'0.006526"
$ws.Cells.Item(42, 5).Value = "'-10.00%"
$ws.Cells.Item(43, 4).Value = "'0.002082"
$ws.Cells.Item(43, 5).Value = "'-1.82%"
$ws.Cells.Item(44, 4).Value = "'0.007508"
$ws.Cells.Item(44, 5).Value = "'-10.47%"
$ws.Cells.Item(45, 4).Value = "'0.3483"
$ws.Cells.Item(45, 5).Value = "'-2.07%"
$ws.Cells.Item(46, 4).Value = "'0.00006973"
$ws.Cells.Item(46, 5).Value = "'-4.34%"
$ws.Cells.Item(47, 4).Value = "'0.00000000749"
$ws.Cells.Item(47, 5).Value = "'-0.45%"
$ws.Cells.Item(48, 4).Value = "'0.003469"
$ws.Cells.Item(48, 5).Value = "'-0.87%"
$ws.Cells.Item(49, 4).Value = "'0.003525"
$ws.Cells.Item(49, 5).Value = "'0.47%"
$ws.Cells.Item(50, 4).Value = "'0.00002097"
$ws.Cells.Item(50, 5).Value = "'-0.45%"
$ws.Cells.Item(51, 4).Value = "'0.0001997"
$ws.Cells.Item(51, 5).Value = "'-0.45%"
